$wb = $excel.ActiveWorkbook

# --- 1. Update the devices sheet values ---
$devices = $wb.Worksheets.Item("devices")
$devices.Range("E2").Value = "Galaxy S7 Edge"
$devices.Range("L2").Value = "au.gov.nsw.onegov.app.checker.uat"

# --- 2. Create "checkerapp" sheet right after "devices" (mirrors devices layout) ---
$checkerapp = $wb.Worksheets.Add($null, $devices)
$checkerapp.Name = "checkerapp"

$checkerapp.Range("A1").Value = "deviceName"
$checkerapp.Range("B1").Value = "platformName"
$checkerapp.Range("C1").Value = "platformVersion"
$checkerapp.Range("D1").Value = "manufacturer"
$checkerapp.Range("E1").Value = "model"
$checkerapp.Range("F1").Value = "resolution"
$checkerapp.Range("G1").Value = "network"
$checkerapp.Range("H1").Value = "location"
$checkerapp.Range("I1").Value = "description"
$checkerapp.Range("J1").Value = "browserName"
$checkerapp.Range("K1").Value = "automationName"
$checkerapp.Range("L1").Value = "appPackage"
$checkerapp.Range("A1:L1").Interior.Color = 65535

$checkerapp.Range("B2").Value = "Android"
$checkerapp.Range("E2").Value = "Galaxy S7 Edge"
$checkerapp.Range("K2").Value = "PerfectoMobile"
$checkerapp.Range("L2").Value = "au.gov.nsw.onegov.app.checker.uat"

$checkerapp.Columns.Item(1).ColumnWidth = 18.592447916666668
$checkerapp.Columns.Item(2).ColumnWidth = 15.166666666666666
$checkerapp.Columns.Item(3).ColumnWidth = 17.022135416666668
$checkerapp.Columns.Item(4).ColumnWidth = 15.877604166666666
$checkerapp.Columns.Item(9).ColumnWidth = 15.877604166666666
$checkerapp.Columns.Item(10).ColumnWidth = 13.166666666666666
$checkerapp.Columns.Item(11).ColumnWidth = 15.307291666666666
$checkerapp.Columns.Item(12).ColumnWidth = 36.451822916666664

$checkerapp.Range("L2").Select() | Out-Null

# --- 3. Create "checkerSignIn" sheet right after "checkerapp" ---
$checkerSignIn = $wb.Worksheets.Add($null, $checkerapp)
$checkerSignIn.Name = "checkerSignIn"

$checkerSignIn.Range("A1").Value = "Username"
$checkerSignIn.Range("B1").Value = "Password"
$checkerSignIn.Range("C1").Value = "pin"
$checkerSignIn.Range("D1").Value = "appName"
$checkerSignIn.Range("A1:D1").Interior.Color = 65535

$checkerSignIn.Range("A2").Value = "DPI-DLTEST"
$checkerSignIn.Range("B2").Value = "SNSW-1234"
$checkerSignIn.Range("C2").Value = "'1234"
$checkerSignIn.Range("D2").Value = "au.gov.nsw.onegov.app.checker.uat"

$checkerSignIn.Columns.Item(1).ColumnWidth = 41.307291666666664
$checkerSignIn.Columns.Item(2).ColumnWidth = 23.307291666666668
$checkerSignIn.Columns.Item(4).ColumnWidth = 32.022135416666664

$checkerSignIn.Range("D1").Select() | Out-Null

# --- 4. Update the devices sheet selection ---
$devices.Range("E8").Select() | Out-Null

# --- 5. Update signIn sheet selection ---
$signIn = $wb.Worksheets.Item("signIn")
$signIn.Range("A2").Select() | Out-Null

# --- 6. Activate checkerSignIn so it becomes the active tab ---
$checkerSignIn.Activate()
